$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update execution time and start/end times for row 2 (TC_001)
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = "2025-03-10 18:36:12"
$ws.Range("I2").Value = "2025-03-10 18:36:12"

# Update start/end times for row 3 (TC_002)
$ws.Range("H3").Value = "2025-03-10 18:36:21"
$ws.Range("I3").Value = "2025-03-10 18:36:22"
